# The row with Temperature = 200 deg C (row 12 on the "Data" sheet) was
# recorded with a clearly erroneous Xylose concentration (T12 = 5221,
# vs. every other value in that column being < 10). Per the commit
# message ("recording errors in SimpleMl") this bad data row is removed
# entirely, so Excel shifts all the rows below it up by one.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet4")
$ws2 = $wb.Worksheets.Item("Data")

# Delete the whole erroneous row (entire-row delete shifts rows 13:26 -> 12:25).
$ws2.Activate()
$ws2.Rows(12).Delete()

# Restore the view/selection state left on the "Sheet4" tab (it was the
# other sheet touched by the session) without leaving it as the active tab.
$ws1.Activate()
$ws1.Range("F28").Select()

# Finish back on "Data" (it stays the active/visible tab), with the
# selection Excel leaves right after a whole-row delete: the emptied
# row is still selected in full, active cell resting in column B.
$ws2.Activate()
$ws2.Range("A12:XFD12").Select()
$ws2.Range("B12").Activate()
